# Macroferia Regional de Talca - Poroto verde
# Insert one new weekly price observation as row 44, pushing the existing
# rows 44-128 down to 45-129 (hence dimension grows from A1:R128 to A1:R129).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 44..128 down by inserting a fresh row at 44.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(44, 1).Value  = 5
$ws.Cells.Item(44, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(44, 3).Value  = "Maule"
$ws.Cells.Item(44, 4).Value  = 44557
$ws.Cells.Item(44, 5).Value  = 7
$ws.Cells.Item(44, 6).Value  = 100112031
$ws.Cells.Item(44, 7).Value  = "Poroto verde"
$ws.Cells.Item(44, 8).Value  = "Sin especificar"
$ws.Cells.Item(44, 9).Value  = "Primera"
$ws.Cells.Item(44, 10).Value = 200
$ws.Cells.Item(44, 11).Value = 40000
$ws.Cells.Item(44, 12).Value = 40000
$ws.Cells.Item(44, 13).Value = 40000
$ws.Cells.Item(44, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(44, 15).Value = "Región del Maule"
$ws.Cells.Item(44, 16).Value = 1600
$ws.Cells.Item(44, 17).Value = 25
$ws.Cells.Item(44, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# rest of column D (style carries over automatically from Insert, but set
# it explicitly to be safe).
$ws.Cells.Item(44, 4).NumberFormat = $ws.Cells.Item(45, 4).NumberFormat
